# "version A simulation rsults"
#
# The author re-saved the workbook after:
#   1. Scrolling/selecting cell G6 on Sheet1 (clearing the old "A30" scroll
#      position and the old K45 selection that were left over from a
#      previous session).
#   2. Removing the color-scale conditional formatting that had been
#      applied over C4:F43.
#
# (Excel's own save routine also renumbers/re-orders the internal style
# table at the same time, but that produces no visible formatting change
# at all -- every cell keeps exactly the same number format, borders and
# fill -- so there is nothing else to replicate here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the color-scale conditional formatting rule over C4:F43.
[void]$ws.Range("C4:F43").FormatConditions.Delete()

# Move the selection to G6; this also resets the window's scroll position,
# so the stale topLeftCell="A30" is dropped.
$ws.Activate() | Out-Null
[void]$ws.Range("G6").Select()
